$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.65
$ws.Range("J2").Value = 4.7
$ws.Range("K2").Value = 4.9
$ws.Range("L2").Value = 1.31
$ws.Range("N2").Value = 5.5
$ws.Range("O2").Value = 1.21
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 1.65
$ws.Range("R2").Value = 1.6
$ws.Range("S2").Value = 2.6
$ws.Range("T2").Value = 1.71
$ws.Range("U2").Value = 2.34
$ws.Range("V2").Value = 1.21
$ws.Range("W2").Value = 2.54
$ws.Range("Y2").Value = 26
$ws.Range("Z2").Value = 50
$ws.Range("AB2").Value = 11.5
$ws.Range("AE2").Value = 65
$ws.Range("AF2").Value = 11.5
$ws.Range("AH2").Value = 18.5
$ws.Range("AK2").Value = 14.5
$ws.Range("AL2").Value = 27
$ws.Range("AN2").Value = 7

# Row 3
$ws.Range("F3").Value = 29
$ws.Range("G3").Value = 38
$ws.Range("I3").Value = 1.11
$ws.Range("J3").Value = 14.5
$ws.Range("K3").Value = 20
$ws.Range("L3").Value = 1.1
$ws.Range("P3").Value = 6.2
$ws.Range("Q3").Value = 1.14
$ws.Range("R3").Value = 3.05
$ws.Range("S3").Value = 1.37
$ws.Range("T3").Value = 1.78
$ws.Range("U3").Value = 1.92
$ws.Range("V3").Value = 9.199999999999999
$ws.Range("X3").Value = 1000
$ws.Range("AA3").Value = 15
$ws.Range("AD3").Value = 23
$ws.Range("AE3").Value = 19
$ws.Range("AH3").Value = 65
$ws.Range("AI3").Value = 44
$ws.Range("AO3").Value = 1.9

# Row 4
$ws.Range("F4").Value = 2.24
$ws.Range("H4").Value = 4.1
$ws.Range("I4").Value = 4.8
$ws.Range("J4").Value = 2.76
$ws.Range("K4").Value = 3.05
$ws.Range("L4").Value = 1.71
$ws.Range("M4").Value = 1.17
$ws.Range("N4").Value = 2.24
$ws.Range("O4").Value = 1.73
$ws.Range("P4").Value = 1.39
$ws.Range("Q4").Value = 3.3
$ws.Range("R4").Value = 1.13
$ws.Range("S4").Value = 7.6
$ws.Range("T4").Value = 2.46
$ws.Range("U4").Value = 1.56
$ws.Range("V4").Value = 1.27
$ws.Range("W4").Value = 1.7
$ws.Range("X4").Value = 6.8
$ws.Range("Y4").Value = 10
$ws.Range("Z4").Value = 110
$ws.Range("AB4").Value = 6.2
$ws.Range("AC4").Value = 7.2
$ws.Range("AF4").Value = 12.5
$ws.Range("AG4").Value = 22
$ws.Range("AH4").Value = 110
$ws.Range("AJ4").Value = 150
$ws.Range("AL4").Value = 460

# Row 5
$ws.Range("F5").Value = 2.26
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 3.1
$ws.Range("K5").Value = 4.3
$ws.Range("L5").Value = 1.27
$ws.Range("N5").Value = 6.4
$ws.Range("P5").Value = 2.8
$ws.Range("Q5").Value = 1.49
$ws.Range("R5").Value = 1.74
$ws.Range("S5").Value = 2.22
$ws.Range("U5").Value = 2.82
$ws.Range("V5").Value = 1.47
$ws.Range("W5").Value = 1.72
$ws.Range("X5").Value = 75
$ws.Range("Y5").Value = 22
$ws.Range("Z5").Value = 65
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 18.5
$ws.Range("AF5").Value = 55
$ws.Range("AH5").Value = 15
$ws.Range("AJ5").Value = 80
$ws.Range("AK5").Value = 22
$ws.Range("AL5").Value = 48
$ws.Range("AN5").Value = 10.5
$ws.Range("AO5").Value = 16.5

# Row 6
$ws.Range("F6").Value = 2.6
$ws.Range("G6").Value = 2.94
$ws.Range("H6").Value = 3.2
$ws.Range("I6").Value = 3.8
$ws.Range("J6").Value = 2.66
$ws.Range("L6").Value = 1.71
$ws.Range("M6").Value = 1.17
$ws.Range("N6").Value = 2.26
$ws.Range("O6").Value = 1.68
$ws.Range("P6").Value = 1.39
$ws.Range("Q6").Value = 3.1
$ws.Range("R6").Value = 1.13
$ws.Range("T6").Value = 1.04
$ws.Range("U6").Value = 1.04
$ws.Range("W6").Value = 1.53
$ws.Range("X6").Value = 1000
$ws.Range("Y6").Value = 1000
$ws.Range("Z6").Value = 980
$ws.Range("AA6").Value = 190
$ws.Range("AB6").Value = 1000
$ws.Range("AC6").Value = 970
$ws.Range("AD6").Value = 1000
$ws.Range("AE6").Value = 450
$ws.Range("AF6").Value = 1000
$ws.Range("AG6").Value = 1000
$ws.Range("AH6").Value = 990
$ws.Range("AI6").Value = 1000
$ws.Range("AJ6").Value = 980
$ws.Range("AK6").Value = 980
$ws.Range("AL6").Value = 1000
$ws.Range("AN6").Value = 1000
$ws.Range("AO6").Value = 1000

# Row 7
$ws.Range("F7").Value = 1.43
$ws.Range("G7").Value = 1.48
$ws.Range("H7").Value = 8
$ws.Range("I7").Value = 9.6
$ws.Range("J7").Value = 4.8
$ws.Range("K7").Value = 5.3
$ws.Range("L7").Value = 1.35
$ws.Range("N7").Value = 4.4
$ws.Range("O7").Value = 1.26
$ws.Range("P7").Value = 2.16
$ws.Range("Q7").Value = 1.77
$ws.Range("S7").Value = 2.92
$ws.Range("T7").Value = 1.95
$ws.Range("U7").Value = 1.86
$ws.Range("V7").Value = 1.11
$ws.Range("W7").Value = 3.05
$ws.Range("Y7").Value = 990
$ws.Range("Z7").Value = 80
$ws.Range("AA7").Value = 340
$ws.Range("AB7").Value = 8.4
$ws.Range("AD7").Value = 34
$ws.Range("AE7").Value = 140
$ws.Range("AF7").Value = 8.800000000000001
$ws.Range("AH7").Value = 27
$ws.Range("AI7").Value = 130
$ws.Range("AM7").Value = 170
$ws.Range("AO7").Value = 170

# Row 8
$ws.Range("G8").Value = 2.88
$ws.Range("H8").Value = 2.72
$ws.Range("L8").Value = 1.44
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 3.6
$ws.Range("O8").Value = 1.36
$ws.Range("P8").Value = 1.86
$ws.Range("Q8").Value = 2.06
$ws.Range("R8").Value = 1.32
$ws.Range("S8").Value = 3.7
$ws.Range("T8").Value = 1.76
$ws.Range("AA8").Value = 900
$ws.Range("AE8").Value = 90
$ws.Range("AN8").Value = 70

# Row 9
$ws.Range("H9").Value = 9
$ws.Range("I9").Value = 11.5
$ws.Range("L9").Value = 1.68
$ws.Range("M9").Value = 1.15
$ws.Range("N9").Value = 2.26
$ws.Range("O9").Value = 1.68
$ws.Range("P9").Value = 1.43
$ws.Range("Q9").Value = 3.1
$ws.Range("R9").Value = 1.14
$ws.Range("S9").Value = 6.6
$ws.Range("T9").Value = 2.78
$ws.Range("Y9").Value = 1000
$ws.Range("AG9").Value = 36

# Row 10
$ws.Range("H10").Value = 7.2
$ws.Range("I10").Value = 7.4
$ws.Range("J10").Value = 3.7
$ws.Range("K10").Value = 3.75
$ws.Range("N10").Value = 3.05
$ws.Range("O10").Value = 1.47
$ws.Range("P10").Value = 1.68
$ws.Range("Q10").Value = 2.44
$ws.Range("T10").Value = 2.32
$ws.Range("U10").Value = 1.72
$ws.Range("W10").Value = 2.44
$ws.Range("X10").Value = 10
$ws.Range("Y10").Value = 17.5
$ws.Range("AA10").Value = 230
$ws.Range("AI10").Value = 150
$ws.Range("AJ10").Value = 15.5
$ws.Range("AM10").Value = 230
$ws.Range("AO10").Value = 220

# Row 11
$ws.Range("F11").Value = 1.12
$ws.Range("G11").Value = 1.13
$ws.Range("H11").Value = 44
$ws.Range("J11").Value = 11.5
$ws.Range("K11").Value = 12
$ws.Range("L11").Value = 1.24
$ws.Range("P11").Value = 3.15
$ws.Range("R11").Value = 1.87
$ws.Range("S11").Value = 2.1
$ws.Range("T11").Value = 2.78
$ws.Range("U11").Value = 1.5
$ws.Range("W11").Value = 9.199999999999999
$ws.Range("X11").Value = 48
$ws.Range("Z11").Value = 580
$ws.Range("AB11").Value = 12.5
$ws.Range("AC11").Value = 30
$ws.Range("AF11").Value = 7.6
$ws.Range("AG11").Value = 17
$ws.Range("AH11").Value = 95
$ws.Range("AK11").Value = 16.5
$ws.Range("AL11").Value = 1000
$ws.Range("AM11").Value = 550
$ws.Range("AN11").Value = 2.92
